# Updated: st 24. 06. 2021
# Apply corrected AgTests (F) / AgPosit (G) figures for several historical rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F321" = 90944
    "F325" = 769894
    "G325" = 6466
    "F327" = 224318
    "G327" = 2718
    "F328" = 180672
    "G328" = 2665
    "F329" = 73331
    "G329" = 1723
    "F341" = 283709
    "F345" = 292750
    "F352" = 306712
    "G352" = 3544
    "F368" = 345890
    "G368" = 2295
    "F376" = 222032
    "F402" = 722027
    "F444" = 100036
    "F445" = 81920
    "F446" = 86597
    "G446" = 264
    "F447" = 64499
    "G447" = 182
    "F454" = 50269
    "F455" = 49989
    "F456" = 47793
    "F460" = 55663
    "F467" = 50148
    "G467" = 73
    "F468" = 40488
    "F469" = 38994
    "G469" = 37
    "F470" = 41300
    "F471" = 61858
    "G471" = 49
    "F472" = 47026
    "F473" = 37908
    "G473" = 38
    "F474" = 44103
    "F475" = 31438
    "G475" = 25
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
